$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all updated cells as Text to preserve the source data format
# (prices/percentages/hours are stored as text strings in this sheet).
$updates = @(
    @{ Ref = "D2"; Value = "246.30" }
    @{ Ref = "E2"; Value = "0.81%" }
    @{ Ref = "G2"; Value = "15" }
    @{ Ref = "D3"; Value = "29.29" }
    @{ Ref = "E3"; Value = "7.31%" }
    @{ Ref = "G3"; Value = "15" }
    @{ Ref = "D4"; Value = "5.193" }
    @{ Ref = "E4"; Value = "2.96%" }
    @{ Ref = "G4"; Value = "15" }
    @{ Ref = "D5"; Value = "0.05702" }
    @{ Ref = "E5"; Value = "0.51%" }
    @{ Ref = "G5"; Value = "15" }
    @{ Ref = "D6"; Value = "6.569" }
    @{ Ref = "E6"; Value = "1.45%" }
    @{ Ref = "G6"; Value = "15" }
    @{ Ref = "B7"; Value = "GateToken" }
    @{ Ref = "C7"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" }
    @{ Ref = "D7"; Value = "3.057" }
    @{ Ref = "E7"; Value = "1.92%" }
    @{ Ref = "G7"; Value = "15" }
    @{ Ref = "B8"; Value = "MXToken" }
    @{ Ref = "C8"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" }
    @{ Ref = "D8"; Value = "0.8579" }
    @{ Ref = "E8"; Value = "4.22%" }
    @{ Ref = "G8"; Value = "15" }
    @{ Ref = "B9"; Value = "FTXToken" }
    @{ Ref = "C9"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" }
    @{ Ref = "D9"; Value = "0.8761" }
    @{ Ref = "E9"; Value = "3.78%" }
    @{ Ref = "G9"; Value = "15" }
    @{ Ref = "B10"; Value = "WazirX" }
    @{ Ref = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" }
    @{ Ref = "D10"; Value = "0.1370" }
    @{ Ref = "E10"; Value = "3.38%" }
    @{ Ref = "G10"; Value = "15" }
    @{ Ref = "B11"; Value = "MandalaExchangeToken" }
    @{ Ref = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Ref = "D11"; Value = "0.07077" }
    @{ Ref = "E11"; Value = "2.28%" }
    @{ Ref = "G11"; Value = "15" }
    @{ Ref = "B12"; Value = "BitrueCoin" }
    @{ Ref = "C12"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Ref = "D12"; Value = "0.02870" }
    @{ Ref = "E12"; Value = "0.27%" }
    @{ Ref = "G12"; Value = "15" }
    @{ Ref = "B13"; Value = "BitMartToken" }
    @{ Ref = "C13"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Ref = "D13"; Value = "0.09391" }
    @{ Ref = "E13"; Value = "0.09%" }
    @{ Ref = "G13"; Value = "15" }
    @{ Ref = "B14"; Value = "BitForexToken" }
    @{ Ref = "C14"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Ref = "D14"; Value = "0.001516" }
    @{ Ref = "E14"; Value = "0.11%" }
    @{ Ref = "G14"; Value = "15" }
    @{ Ref = "B15"; Value = "CoinExToken" }
    @{ Ref = "C15"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" }
    @{ Ref = "D15"; Value = "0.04160" }
    @{ Ref = "E15"; Value = "1.06%" }
    @{ Ref = "G15"; Value = "15" }
    @{ Ref = "B16"; Value = "One" }
    @{ Ref = "C16"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" }
    @{ Ref = "D16"; Value = "0.0006021" }
    @{ Ref = "E16"; Value = "-93.99%" }
    @{ Ref = "G16"; Value = "15" }
    @{ Ref = "B17"; Value = "TigerCash" }
    @{ Ref = "C17"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Ref = "D17"; Value = "0.006062" }
    @{ Ref = "E17"; Value = "-1.01%" }
    @{ Ref = "G17"; Value = "15" }
    @{ Ref = "B18"; Value = "UpBots" }
    @{ Ref = "C18"; Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt" }
    @{ Ref = "D18"; Value = "0.007491" }
    @{ Ref = "E18"; Value = "5,107.82%" }
    @{ Ref = "G18"; Value = "15" }
    @{ Ref = "B19"; Value = "LEO" }
    @{ Ref = "C19"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Ref = "D19"; Value = "3.483" }
    @{ Ref = "E19"; Value = "-0.77%" }
    @{ Ref = "G19"; Value = "15" }
    @{ Ref = "D20"; Value = "2.277" }
    @{ Ref = "E20"; Value = "-1.42%" }
    @{ Ref = "G20"; Value = "15" }
    @{ Ref = "D21"; Value = "0.3143" }
    @{ Ref = "E21"; Value = "0.94%" }
    @{ Ref = "G21"; Value = "15" }
    @{ Ref = "D22"; Value = "0.03290" }
    @{ Ref = "E22"; Value = "3.27%" }
    @{ Ref = "G22"; Value = "15" }
    @{ Ref = "D23"; Value = "0.1300" }
    @{ Ref = "E23"; Value = "0.71%" }
    @{ Ref = "G23"; Value = "15" }
    @{ Ref = "D24"; Value = "3.469" }
    @{ Ref = "E24"; Value = "-2.57%" }
    @{ Ref = "G24"; Value = "15" }
    @{ Ref = "G25"; Value = "15" }
    @{ Ref = "D26"; Value = "0.005087" }
    @{ Ref = "E26"; Value = "31.46%" }
    @{ Ref = "G26"; Value = "15" }
    @{ Ref = "D27"; Value = "0.001219" }
    @{ Ref = "E27"; Value = "0.11%" }
    @{ Ref = "G27"; Value = "15" }
    @{ Ref = "E28"; Value = "23.46%" }
    @{ Ref = "G28"; Value = "15" }
    @{ Ref = "G29"; Value = "15" }
    @{ Ref = "G30"; Value = "15" }
    @{ Ref = "G31"; Value = "15" }
    @{ Ref = "G32"; Value = "15" }
    @{ Ref = "G33"; Value = "15" }
    @{ Ref = "G34"; Value = "15" }
    @{ Ref = "G35"; Value = "15" }
    @{ Ref = "G36"; Value = "15" }
    @{ Ref = "G37"; Value = "15" }
    @{ Ref = "G38"; Value = "15" }
    @{ Ref = "G39"; Value = "15" }
    @{ Ref = "D40"; Value = "0.03745" }
    @{ Ref = "E40"; Value = "2.21%" }
    @{ Ref = "G40"; Value = "15" }
    @{ Ref = "D41"; Value = "0.005719" }
    @{ Ref = "E41"; Value = "-5.31%" }
    @{ Ref = "G41"; Value = "15" }
    @{ Ref = "D42"; Value = "0.1069" }
    @{ Ref = "E42"; Value = "1.59%" }
    @{ Ref = "G42"; Value = "15" }
    @{ Ref = "D43"; Value = "0.002100" }
    @{ Ref = "E43"; Value = "-8.41%" }
    @{ Ref = "G43"; Value = "15" }
    @{ Ref = "D44"; Value = "0.009398" }
    @{ Ref = "E44"; Value = "-3.14%" }
    @{ Ref = "G44"; Value = "15" }
    @{ Ref = "D45"; Value = "0.00005109" }
    @{ Ref = "E45"; Value = "-3.95%" }
    @{ Ref = "G45"; Value = "15" }
    @{ Ref = "E46"; Value = "-0.01%" }
    @{ Ref = "G46"; Value = "15" }
    @{ Ref = "G47"; Value = "15" }
    @{ Ref = "D48"; Value = "0.002661" }
    @{ Ref = "E48"; Value = "3.74%" }
    @{ Ref = "G48"; Value = "15" }
    @{ Ref = "E49"; Value = "-0.01%" }
    @{ Ref = "G49"; Value = "15" }
    @{ Ref = "E50"; Value = "-0.01%" }
    @{ Ref = "G50"; Value = "15" }
    @{ Ref = "G51"; Value = "15" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
